$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column V: a baseline "aggregation" flag/header block, appended after
# column U. Shared-string insertion order matters (AGG[95:96], 0 = NOT,
# 1 = AGGREGATING), so write V1 then V3 then V2.
$ws.Range("V1").Value = "AGG[95:96]"
$ws.Range("V3").Value = "0 = NOT"
$ws.Range("V2").Value = "1 = AGGREGATING"

# Widen column V a bit for the new header text, and nudge the view/selection
# to roughly match where the author had scrolled/selected afterwards.
$ws.Columns.Item(22).ColumnWidth = 15.3

$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 12
$ws.Range("U5").Select()
